$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 9
$ws.Range("B2").Value = "SORO GLICOSADO 5%, BOLSA SISTEMA FECHADO, FRASCO COM 1000ML"
$ws.Range("C2").Value = "JP INDUSTRIA FARMACEUTICA S/A"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "104910020"
$ws.Range("E2").Value = "50 MG/ML SOL INJ IV CX 10 BOLS PLAS SIST FECH X 1000 ML Ativo"
$ws.Range("F2").Value = "OK"

$wb.Save()
